# This script applies the edit described by the commit:
# "add number_of_shares in security.py and excel_writer.py"
#
# For each sheet (Alternative, Bond, Equity) it:
#  1. Adds a new "Number of Shares" column in column S, with the same
#     header style as the other header cells (copied from column R).
#  2. Fills in the per-row "Number of Shares" values.
#  3. Updates the recalculated metrics (Simple/Total Return, Std Dev,
#     Downside Deviation, Value at Risk, Sharpe Ratio, Portfolio Asset
#     Weight/Allocation, etc.) that shifted as a result of the new data.

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------
# Sheet: Alternative
# ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Alternative")

# Add the new "Number of Shares" header in S1, copying the style
# (bold, bordered, centered) used by the other header cells.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$ws.Range("S1").Value = "Number of Shares"

# Updated / new numeric values
$ws.Range("K2").Value = 9.31
$ws.Range("L2").Value = 9.31
$ws.Range("M2").Value = 12.19
$ws.Range("O2").Value = -10.61
$ws.Range("P2").Value = 0.6
$ws.Range("S2").Value = 11.73
$ws.Range("K3").Value = 2.98
$ws.Range("L3").Value = 2.98
$ws.Range("M3").Value = 14.79
$ws.Range("N3").Value = 8.699999999999999
$ws.Range("O3").Value = -21.46
$ws.Range("P3").Value = 0.07000000000000001
$ws.Range("S3").Value = 5.16
$ws.Range("K4").Value = 1.25
$ws.Range("L4").Value = 9.5
$ws.Range("M4").Value = 8.75
$ws.Range("N4").Value = 5.33
$ws.Range("O4").Value = -4.71
$ws.Range("S4").Value = 78.89
$ws.Range("K5").Value = 1.52
$ws.Range("L5").Value = 3.38
$ws.Range("M5").Value = 7.76
$ws.Range("P5").Value = 0.18
$ws.Range("S5").Value = 0
$ws.Range("K6").Value = -9.08
$ws.Range("L6").Value = -9.08
$ws.Range("M6").Value = 44.58
$ws.Range("N6").Value = 24.83
$ws.Range("O6").Value = -83.44
$ws.Range("S6").Value = 0
$ws.Range("K7").Value = -0.32
$ws.Range("L7").Value = 2.15
$ws.Range("M7").Value = 6.46
$ws.Range("O7").Value = -8.48
$ws.Range("P7").Value = 0.02
$ws.Range("S7").Value = 0
$ws.Range("K8").Value = 1.57
$ws.Range("L8").Value = 1.58
$ws.Range("M8").Value = 5.48
$ws.Range("O8").Value = -7.27
$ws.Range("P8").Value = -0.08
$ws.Range("S8").Value = 0

# ----------------------------------------------------------------
# Sheet: Bond
# ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Bond")

# Add the new "Number of Shares" header in S1, copying the style
# (bold, bordered, centered) used by the other header cells.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$ws.Range("S1").Value = "Number of Shares"

# Updated / new numeric values
$ws.Range("K2").Value = -3.6
$ws.Range("L2").Value = 10.68
$ws.Range("M2").Value = 10.94
$ws.Range("O2").Value = -7.26
$ws.Range("P2").Value = 0.79
$ws.Range("Q2").Value = 15.75
$ws.Range("R2").Value = 1575
$ws.Range("S2").Value = 78.28
$ws.Range("J3").Value = 3.08
$ws.Range("K3").Value = -0.4
$ws.Range("L3").Value = 4.33
$ws.Range("M3").Value = 0.32
$ws.Range("N3").Value = 2.41
$ws.Range("O3").Value = 3.81
$ws.Range("P3").Value = 7.4
$ws.Range("Q3").Value = 2.63
$ws.Range("R3").Value = 262.64
$ws.Range("S3").Value = 2.66
$ws.Range("H4").Value = 0.07000000000000001
$ws.Range("J4").Value = 3.2
$ws.Range("K4").Value = -1.26
$ws.Range("L4").Value = 3.74
$ws.Range("M4").Value = 18.11
$ws.Range("O4").Value = -26.33
$ws.Range("Q4").Value = 2.63
$ws.Range("R4").Value = 262.64
$ws.Range("S4").Value = 3.28
$ws.Range("O5").Value = -16.85
$ws.Range("Q5").Value = 7
$ws.Range("R5").Value = 700
$ws.Range("S5").Value = 10.98
$ws.Range("K6").Value = 20.41
$ws.Range("L6").Value = 21.12
$ws.Range("M6").Value = 17.65
$ws.Range("O6").Value = -8.07
$ws.Range("Q6").Value = 12
$ws.Range("R6").Value = 1200
$ws.Range("S6").Value = 35.22

# ----------------------------------------------------------------
# Sheet: Equity
# ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("Equity")

# Add the new "Number of Shares" header in S1, copying the style
# (bold, bordered, centered) used by the other header cells.
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)
$ws.Range("S1").Value = "Number of Shares"

# Updated / new numeric values
$ws.Range("H2").Value = 0.14
$ws.Range("J2").Value = 1.48
$ws.Range("K2").Value = 13.64
$ws.Range("L2").Value = 15.14
$ws.Range("M2").Value = 20.11
$ws.Range("O2").Value = -17.79
$ws.Range("Q2").Value = 4.8
$ws.Range("R2").Value = 479.52
$ws.Range("S2").Value = 1.1
$ws.Range("K3").Value = 8.880000000000001
$ws.Range("L3").Value = 11.78
$ws.Range("M3").Value = 16.7
$ws.Range("O3").Value = -16.23
$ws.Range("P3").Value = 0.59
$ws.Range("Q3").Value = 2.1
$ws.Range("R3").Value = 209.88
$ws.Range("S3").Value = 6.51
$ws.Range("K4").Value = 3.86
$ws.Range("L4").Value = 4.7
$ws.Range("M4").Value = 16.39
$ws.Range("O4").Value = -21.94
$ws.Range("Q4").Value = 2.9
$ws.Range("R4").Value = 289.98
$ws.Range("S4").Value = 11.31
$ws.Range("J5").Value = 2.45
$ws.Range("K5").Value = 5.75
$ws.Range("L5").Value = 8.81
$ws.Range("M5").Value = 15.59
$ws.Range("O5").Value = -16.71
$ws.Range("Q5").Value = 1.8
$ws.Range("R5").Value = 180
$ws.Range("S5").Value = 4.55
$ws.Range("K6").Value = 4.55
$ws.Range("L6").Value = 8.869999999999999
$ws.Range("O6").Value = -8.76
$ws.Range("P6").Value = 0.65
$ws.Range("Q6").Value = 2.5
$ws.Range("R6").Value = 250.02
$ws.Range("S6").Value = 9.07
$ws.Range("J7").Value = 2.54
$ws.Range("K7").Value = 6.22
$ws.Range("L7").Value = 8.6
$ws.Range("M7").Value = 13.72
$ws.Range("O7").Value = -14.03
$ws.Range("P7").Value = 0.48
$ws.Range("Q7").Value = 3.91
$ws.Range("R7").Value = 390.78
$ws.Range("S7").Value = 7.05
$ws.Range("J8").Value = 11.74
$ws.Range("K8").Value = -2.23
$ws.Range("L8").Value = 10.71
$ws.Range("M8").Value = 8.710000000000001
$ws.Range("N8").Value = 7.87
$ws.Range("O8").Value = -3.2
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 12
$ws.Range("R8").Value = 1200
$ws.Range("S8").Value = 89.15000000000001
$ws.Range("J9").Value = 2.08
$ws.Range("K9").Value = 11.17
$ws.Range("L9").Value = 11.82
$ws.Range("M9").Value = 31.33
$ws.Range("N9").Value = 10.79
$ws.Range("O9").Value = -38.94
$ws.Range("P9").Value = 0.31
$ws.Range("S9").Value = 0

